$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DDF")

# New row of data (row 6): abc | xyz | '111 (quoted/text) | dispatch message
$ws.Range("A6").Value = "abc"
$ws.Range("B6").Value = "xyz"
# Insert the dispatch message first so it lands at shared-string index 40,
# then the quoted "111" text lands at index 41 (matches source order).
$ws.Range("D6").Value = "Your order has been dispatched, and will arrive just as fast as the pony can get there!"
$ws.Range("C6").Formula = "'111"

# Give new column D a wide custom width (text column holding the long message)
$ws.Columns.Item(4).ColumnWidth = 71.83333333333333

# Move/restore the active cell selection as recorded in the saved view
$ws.Range("D10").Select()
